$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (06/04/2020 per serial date 43986) appended below the
# existing table data (previously rows 1:83, table Condicion_Pacientes).
# Copy formatting from the last existing data row (83) so the new row
# picks up the same date / centered-number styles, then overwrite values.
$ws.Range("A83:F83").Copy()
$ws.Range("A84:F84").PasteSpecial(-4122)

$ws.Range("A84").Value = 43986
$ws.Range("B84").Value = 524
$ws.Range("C84").Value = 190
$ws.Range("D84").Value = 371
$ws.Range("E84").Value = 134
$ws.Range("F84").Value = 34

# "Clinicamente Estables" (D) gets its own wrap-text + centered style.
$ws.Range("D84").WrapText = $true

# Grow the table (ListObject) so its range/autoFilter cover the new row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F84"))

# Match the saved selection state after the edit.
$ws.Range("E85").Select() | Out-Null
